$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "61.982.45"
$ws.Range("E2").Value = "  -3.46%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "2.996.16"
$ws.Range("E3").Value = "  -4.26%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = "  -0.01%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.56"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -4.69%  "

# Row 6: 'Solana' -> 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.93"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -6.14%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("E7").Value = "  -0.03%  "

# Row 8: 'XRP' -> 'XRP'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range("D9").Value = "3.012.59"
$ws.Range("E9").Value = "  -4.23%  "

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E10").Value = "  -2.96%  "

# Row 11: 'Toncoin' -> 'Toncoin'
$ws.Range("E11").Value = "  -6.69%  "

# Row 12: 'Cardano' -> 'Cardano'
$ws.Range("E12").Value = "  -2.83%  "

# Row 13: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D13").Value = "3.523.38"
$ws.Range("E13").Value = "  -4.02%  "

# Row 14: 'TRON' -> 'TRON'
$ws.Range("E14").Value = "  -1.18%  "

# Row 15: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D15").Value = "61.979.41"
$ws.Range("E15").Value = "  -3.54%  "

# Row 16: 'Avalanche' -> 'Avalanche'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.03"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.06%  "

# Row 17: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D17").Value = "3.003.16"
$ws.Range("E17").Value = "  -3.98%  "

# Row 18: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000148"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -4.49%  "

# Row 19: 'Polkadot' -> 'Polkadot'
$ws.Range("E19").Value = "  -0.58%  "

# Row 20: 'Chainlink' -> 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.09"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.94%  "

# Row 21: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.58"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -6.71%  "

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.77"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.61%  "

# Row 23: 'Dai' -> 'Dai'
$ws.Range("E23").Value = "  -0.01%  "

# Row 24: 'Litecoin' -> 'Litecoin'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.12"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.46%  "

# Row 25: 'Polygon' -> 'WrappedeETH'
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.118.23"
$ws.Range("E25").Value = "  -4.65%  "

# Row 26: 'WrappedeETH' -> 'Polygon'
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.471"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.93%  "

# Row 27: 'Kaspa' -> 'Kaspa'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.191"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.24%  "

# Row 28: 'Binance-PegBSC-USD' -> 'PEPE'
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0949"
$ws.Range("E28").Value = "  -6.59%  "

# Row 29: 'PEPE' -> 'Binance-PegBSC-USD'
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -6.45%  "

# Row 31: 'USDe' -> 'USDe'
$ws.Range("E31").Value = "  -0.02%  "

# Row 32: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.73"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -3.69%  "

# Row 33: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.56"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -3.06%  "

# Row 34: 'Monero' -> 'Monero'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "160.37"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.75%  "

# Row 35: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.64"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -3.47%  "

# Row 36: 'Aptos' -> 'Aptos'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -3.54%  "

# Row 37: 'Fetch.AI' -> 'Fetch.AI'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.08"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.50%  "

# Row 38: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.29"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -4.36%  "

# Row 39: 'Stacks' -> 'Stacks'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -6.77%  "

# Row 40: 'OKB' -> 'OKB'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.73"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.15%  "

# Row 41: 'Filecoin' -> 'Maker'
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.428.67"
$ws.Range("E41").Value = "  -7.50%  "

# Row 42: 'Maker' -> 'Filecoin'
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.92"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -3.93%  "

# Row 43: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.30"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.51%  "

# Row 44: 'Mantle' -> 'Mantle'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.676"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.86%  "

# Row 45: 'Hedera' -> 'Hedera'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0596"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -2.50%  "

# Row 46: 'RenderToken' -> 'RenderToken'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.21"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -1.46%  "

# Row 47: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range("E47").Value = "  -0.06%  "

# Row 48: 'VeChain' -> 'VeChain'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0247"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -2.42%  "

# Row 49: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.98"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -5.75%  "

# Row 50: 'Stellar' -> 'Stellar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0957"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.47%  "

# Row 51: 'Bittensor' -> 'Bittensor'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "270.02"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -6.58%  "
